# "Add files via upload" - corrections made to the product sheet before
# re-uploading the workbook:
#   - D4: fix stray space in the image filename ("foot .jpg" -> "foot.jpg")
#   - F6:F9: capitalize "Amazon choice" -> "Amazon Choice"
#   - Rows 6-9 resized taller while reviewing the sheet
#   - Scrolled/selected down to F9 before saving

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the filename typo
$ws.Range("D4").Value = "foot.jpg"

# Fix capitalization of the "Amazon Choice" label for the four affected rows
$ws.Range("F6:F9").Value = "Amazon Choice"

# Manually resized rows 6-9 (taller to show full wrapped text)
$ws.Rows.Item(6).RowHeight = 105
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45

# Scroll/selection state at save time
$ws.Range("F9").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
